$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 4 (a7e0052c... file):
#   Column E = Correspond Handoff Datetime
#   Column H = Correspond Handback DateTime
$wsZhCn.Range("E4").Value = "2016-03-13 00:37:47"
$wsZhCn.Range("H4").Value = "2016-03-13 00:38:18"

# de-de sheet, row 4 (a7e0052c... file):
$wsDeDe.Range("E4").Value = "2016-03-13 00:37:51"
$wsDeDe.Range("H4").Value = "2016-03-13 00:38:24"
